$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply each updated cell as literal text, preserving the General/
# unstyled cell format (matches the source workbook, which stores every
# data cell as an inline string with no explicit number format).
function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextCell 2 4 '25.869.68'
Set-TextCell 2 5 '  +0.30%  '
Set-TextCell 3 4 '1.636.76'
Set-TextCell 3 5 '  +0.64%  '
Set-TextCell 4 5 '  +0.16%  '
Set-TextCell 5 4 '215.13'
Set-TextCell 5 5 '  -0.01%  '
Set-TextCell 6 4 '0.5091'
Set-TextCell 6 5 '  -0.40%  '
Set-TextCell 7 5 '  +0.22%  '
Set-TextCell 8 4 '0.2585'
Set-TextCell 8 5 '  +0.80%  '
Set-TextCell 9 4 '0.06432'
Set-TextCell 9 5 '  +1.70%  '
Set-TextCell 10 5 '  +4.88%  '
Set-TextCell 11 4 '0.07790'
Set-TextCell 11 5 '  +0.21%  '
Set-TextCell 12 2 'WrappedEther'
Set-TextCell 12 3 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell 12 4 '1.664.43'
Set-TextCell 12 5 '  +1.95%  '
Set-TextCell 13 2 'Polkadot'
Set-TextCell 13 3 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell 13 4 '4.274'
Set-TextCell 13 5 '  +1.10%  '
Set-TextCell 14 4 '1.861.59'
Set-TextCell 14 5 '  +0.70%  '
Set-TextCell 15 4 '0.5607'
Set-TextCell 15 5 '  +1.10%  '
Set-TextCell 16 4 '0.0₅7675'
Set-TextCell 16 5 '  +2.35%  '
Set-TextCell 17 4 '63.26'
Set-TextCell 17 5 '  -0.39%  '
Set-TextCell 18 4 '25.876.68'
Set-TextCell 18 5 '  +0.33%  '
Set-TextCell 19 5 '  +0.28%  '
Set-TextCell 20 4 '193.62'
Set-TextCell 20 5 '  -0.27%  '
Set-TextCell 21 4 '4.387'
Set-TextCell 21 5 '  -0.68%  '
Set-TextCell 22 4 '9.953'
Set-TextCell 22 5 '  +1.95%  '
Set-TextCell 23 4 '6.168'
Set-TextCell 23 5 '  +2.74%  '
Set-TextCell 24 4 '1.004'
Set-TextCell 24 5 '  +0.20%  '
Set-TextCell 25 4 '1.789'
Set-TextCell 26 4 '138.73'
Set-TextCell 26 5 '  -2.01%  '
Set-TextCell 27 5 '  -1.62%  '
Set-TextCell 28 4 '6.861'
Set-TextCell 28 5 '  +2.31%  '
Set-TextCell 29 4 '15.54'
Set-TextCell 29 5 '  +0.26%  '
Set-TextCell 30 4 '1.241'
Set-TextCell 30 5 '  +0.28%  '
Set-TextCell 31 4 '0.04960'
Set-TextCell 31 5 '  +2.30%  '
Set-TextCell 32 5 '  +1.84%  '
Set-TextCell 33 4 '3.253'
Set-TextCell 33 5 '  +2.68%  '
Set-TextCell 34 4 '1.568'
Set-TextCell 34 5 '  +1.98%  '
Set-TextCell 35 5 '  +0.93%  '
Set-TextCell 36 4 '0.9044'
Set-TextCell 36 5 '  +1.32%  '
Set-TextCell 37 5 '  +1.42%  '
Set-TextCell 38 4 '0.5570'
Set-TextCell 38 5 '  +1.19%  '
Set-TextCell 39 4 '1.134.98'
Set-TextCell 39 5 '  +1.94%  '
Set-TextCell 40 4 '0.01569'
Set-TextCell 40 5 '  +1.43%  '
Set-TextCell 41 4 '0.9973'
Set-TextCell 41 5 '  -0.37%  '
Set-TextCell 42 2 'FraxShare'
Set-TextCell 42 3 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell 42 4 '5.469'
Set-TextCell 42 5 '  -0.95%  '
Set-TextCell 43 2 'Quant'
Set-TextCell 43 3 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell 43 4 '99.27'
Set-TextCell 43 5 '  +2.20%  '
Set-TextCell 44 4 '0.8009'
Set-TextCell 44 5 '  +0.69%  '
Set-TextCell 45 5 '  -4.01%  '
Set-TextCell 46 4 '55.51'
Set-TextCell 46 5 '  +1.70%  '
Set-TextCell 47 4 '0.4264'
Set-TextCell 47 5 '  -3.59%  '
Set-TextCell 48 4 '7.746'
Set-TextCell 48 5 '  +2.94%  '
Set-TextCell 49 4 '0.05074'
Set-TextCell 49 5 '  -1.05%  '
Set-TextCell 50 4 '1.001'
Set-TextCell 50 5 '  +0.37%  '
Set-TextCell 51 5 '  +0.33%  '
